$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Change 1: turn the first paragraph
#   "This is a Microsoft word document."
# into
#   "This is a Microsoft word document.  (This is a change - Version
#    for branch alternate)"
# with the parenthetical part colored dark red (C00000), built up as
# three separate insertions (mirrors how the real edit was typed/pasted).
# ----------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r = $p1.Range
[void]$r.MoveEnd(1, -1)    # exclude the paragraph mark
$r.Collapse(0)             # collapse to just after "document."

$r.InsertAfter("  ")
$r.Collapse(0)

$r.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r.Font.Color = 192        # RGB(192,0,0) == C00000
$r.Collapse(0)

$r.InsertAfter("rsion for branch alternate")
$r.Font.Color = 192
$r.Collapse(0)

$r.InsertAfter(")")
$r.Font.Color = 192
$r.Collapse(0)

# ----------------------------------------------------------------------
# Change 2: add a new, empty, shaded (fill F9F9F9) paragraph at the very
# end of the document (after the last "...free at last." paragraph).
# ----------------------------------------------------------------------
[void]$d.Content.Find.Execute("we are free at last.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "we are free at last.^p", 2)

$lastPara = $d.Paragraphs.Last
$lastPara.Style = "Normal"
$lastPara.Shading.Texture = 0
$lastPara.Shading.ForegroundPatternColor = -16777216
$lastPara.Shading.BackgroundPatternColor = 16382457

# ----------------------------------------------------------------------
# Change 3: drop the now-unused custom/heading styles from styles.xml
# (Heading2, Heading4 and their linked character styles, plus a handful
# of leftover custom styles from earlier pasted web content). Deleted in
# reverse definition order so earlier lookups-by-name stay valid.
# ----------------------------------------------------------------------
$staleStyles = @(
  "podcast-tools__subscribe-links",
  "generic-title",
  "subscribe-more-info",
  "subscribe",
  "audio-tool",
  "Heading 4 Char",
  "Heading 2 Char",
  "Hyperlink",
  "apple-converted-space",
  "Heading 4",
  "Heading 2"
)
foreach ($styleName in $staleStyles) {
  $d.Styles($styleName).Delete()
}
